# The post formerly listed at row 673 ("「スマホが好きな理由は…」") was removed
# from the source data. Deleting the entire worksheet row shifts every
# subsequent row up by one (674->673, 675->674, ... 846->845), which matches
# the reference diff (including the updated dimension A1:C846 -> A1:C845).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(673).Delete()
